# Update "想去人数" (interest count, column F) values for a handful of
# events that appear on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# Row -> new value, for the "展览" sheet
$expoUpdates = @{
    3  = 3843
    6  = 3836
    10 = 8684
    11 = 491
    19 = 10997
    21 = 144
    42 = 348
    46 = 130
}

foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Row -> new value, for the "全部类型" sheet
$allUpdates = @{
    3  = 3843
    7  = 3836
    11 = 8684
    12 = 491
    18 = 10997
    21 = 144
    43 = 348
    46 = 130
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
